$d = $word.ActiveDocument

$pairs = @(
    @("676÷6=", "326÷2="),
    @("256÷3=", "469÷9="),
    @("602÷6=", "777÷6="),
    @("200÷3=", "439÷6="),
    @("462÷8=", "692÷5="),
    @("253÷9=", "146÷8="),
    @("452÷9=", "492÷6="),
    @("630÷4=", "835÷4="),
    @("520÷5=", "109÷2="),
    @("824÷5=", "692÷2="),
    @("139÷7=", "188÷9="),
    @("430÷8=", "169÷9="),
    @("333÷7=", "359÷5="),
    @("286÷8=", "329÷6="),
    @("949÷9=", "990÷3="),
    @("441÷4=", "220÷4="),
    @("859÷4=", "935÷2="),
    @("160÷7=", "819÷8="),
    @("648÷9=", "479÷4="),
    @("244÷4=", "577÷5="),
    @("500÷3=", "434÷8="),
    @("484÷4=", "322÷6="),
    @("656÷7=", "869÷3="),
    @("589÷2=", "269÷7="),
    @("479÷5=", "141÷4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
